$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete data: drop the last forecast-origin column (BA) ---
$ws.Range("BA1:BA24").Clear()

# --- Remove obsolete rows 23 and 24 (forecast origins no longer produced) ---
$ws.Range("A23:AZ24").Clear()

# --- Remove stale warm-up forecast values that the bugfix drops from each row ---
$ws.Range("C3:D3").Clear()
$ws.Range("C4:F4").Clear()
$ws.Range("C5:H5").Clear()
$ws.Range("E6:J6").Clear()
$ws.Range("G7:J7").Clear()
$ws.Range("I8:J8").Clear()
$ws.Range("K9").Clear()
$ws.Range("M10").Clear()
$ws.Range("O11").Clear()
$ws.Range("Q12").Clear()
$ws.Range("R13:S13").Clear()
$ws.Range("T14:V14").Clear()
$ws.Range("V15:Z15").Clear()
$ws.Range("Y16:AD16").Clear()
$ws.Range("AC17:AG17").Clear()
$ws.Range("AG18:AK18").Clear()
$ws.Range("AK19:AO19").Clear()
$ws.Range("AO20:AS20").Clear()
$ws.Range("AS21:AW21").Clear()
$ws.Range("AW22:AZ22").Clear()

# --- Update values that the corrected naive-forecaster recomputation produced ---
$ws.Cells.Item(1, 4).Value = 39948
$ws.Cells.Item(1, 2).Value = 39583
$ws.Cells.Item(1, 35).Value = 44341
$ws.Cells.Item(1, 52).Value = 45891
$ws.Cells.Item(1, 15).Value = 41957
$ws.Cells.Item(1, 29).Value = 43783
$ws.Cells.Item(1, 43).Value = 45071
$ws.Cells.Item(1, 50).Value = 45713
$ws.Cells.Item(1, 16).Value = 42137
$ws.Cells.Item(1, 32).Value = 44068
$ws.Cells.Item(1, 49).Value = 45618
$ws.Cells.Item(1, 45).Value = 45254
$ws.Cells.Item(1, 51).Value = 45800
$ws.Cells.Item(1, 19).Value = 42689
$ws.Cells.Item(1, 20).Value = 42867
$ws.Cells.Item(1, 13).Value = 41592
$ws.Cells.Item(1, 25).Value = 43418
$ws.Cells.Item(1, 8).Value = 40676
$ws.Cells.Item(1, 23).Value = 43235
$ws.Cells.Item(1, 34).Value = 44251
$ws.Cells.Item(1, 44).Value = 45163
$ws.Cells.Item(1, 12).Value = 41409
$ws.Cells.Item(1, 36).Value = 44432
$ws.Cells.Item(1, 27).Value = 43600
$ws.Cells.Item(1, 17).Value = 42321
$ws.Cells.Item(1, 5).Value = 40130
$ws.Cells.Item(1, 40).Value = 44798
$ws.Cells.Item(1, 47).Value = 45436
$ws.Cells.Item(1, 38).Value = 44617
$ws.Cells.Item(1, 24).Value = 43326
$ws.Cells.Item(1, 46).Value = 45345
$ws.Cells.Item(1, 11).Value = 41228
$ws.Cells.Item(1, 30).Value = 43875
$ws.Cells.Item(1, 31).Value = 43966
$ws.Cells.Item(1, 33).Value = 44159
$ws.Cells.Item(1, 3).Value = 39765
$ws.Cells.Item(1, 37).Value = 44525
$ws.Cells.Item(1, 48).Value = 45534
$ws.Cells.Item(1, 10).Value = 41044
$ws.Cells.Item(1, 22).Value = 43145
$ws.Cells.Item(1, 28).Value = 43691
$ws.Cells.Item(1, 26).Value = 43510
$ws.Cells.Item(1, 39).Value = 44706
$ws.Cells.Item(1, 7).Value = 40494
$ws.Cells.Item(1, 41).Value = 44890
$ws.Cells.Item(1, 14).Value = 41774
$ws.Cells.Item(1, 42).Value = 44981
$ws.Cells.Item(1, 9).Value = 40862
$ws.Cells.Item(1, 6).Value = 40310
$ws.Cells.Item(1, 21).Value = 43053
$ws.Cells.Item(1, 18).Value = 42503
$ws.Cells.Item(7, 11).Value = -0.3496173419443749
$ws.Cells.Item(7, 12).Value = -0.6616365666142765
$ws.Cells.Item(8, 13).Value = 0.05500386022236903
$ws.Cells.Item(8, 14).Value = 0.2184978785563896
$ws.Cells.Item(8, 11).Value = 0.1626813993622633
$ws.Cells.Item(8, 12).Value = 0.2197847717222867
$ws.Cells.Item(9, 13).Value = 0.05909805310246874
$ws.Cells.Item(9, 12).Value = 0.09544618126309246
$ws.Cells.Item(9, 14).Value = 0.04195831742983547
$ws.Cells.Item(9, 15).Value = 0.07916875696107883
$ws.Cells.Item(9, 16).Value = -0.01790997771649039
$ws.Cells.Item(10, 16).Value = 0.1671491311400208
$ws.Cells.Item(10, 17).Value = 0.1656566557188155
$ws.Cells.Item(10, 15).Value = 0.09103564879091586
$ws.Cells.Item(10, 18).Value = -0.08711135105702317
$ws.Cells.Item(10, 14).Value = 0.08086314912676418
$ws.Cells.Item(11, 19).Value = 0.2044493994367125
$ws.Cells.Item(11, 18).Value = 0.1158714888162216
$ws.Cells.Item(11, 17).Value = 0.09512483792448734
$ws.Cells.Item(11, 20).Value = 0.1761917659537371
$ws.Cells.Item(11, 16).Value = 0.09001325883963851
$ws.Cells.Item(12, 20).Value = 0.1078587431702305
$ws.Cells.Item(12, 24).Value = 0.1415392254179304
$ws.Cells.Item(12, 19).Value = 0.08337410701473313
$ws.Cells.Item(12, 22).Value = 0.2803378563356329
$ws.Cells.Item(12, 18).Value = 0.07295307304728826
$ws.Cells.Item(12, 21).Value = 0.2685680645708288
$ws.Cells.Item(12, 23).Value = 0.3524405906205841
$ws.Cells.Item(13, 22).Value = 0.1089201880626334
$ws.Cells.Item(13, 26).Value = 0.3424613118119479
$ws.Cells.Item(13, 27).Value = 0.07560805834034845
$ws.Cells.Item(13, 21).Value = 0.1068584531153549
$ws.Cells.Item(13, 23).Value = 0.1531204771924033
$ws.Cells.Item(13, 20).Value = 0.08648845420498041
$ws.Cells.Item(13, 25).Value = 0.2977174885593792
$ws.Cells.Item(13, 24).Value = -0.1904239862803969
$ws.Cells.Item(13, 28).Value = -0.05069288950212414
$ws.Cells.Item(14, 30).Value = 0.01105513701109562
$ws.Cells.Item(14, 26).Value = 0.1339497680586277
$ws.Cells.Item(14, 24).Value = 0.0862183575680131
$ws.Cells.Item(14, 32).Value = -2.657403949513992
$ws.Cells.Item(14, 29).Value = 0.0103609600907939
$ws.Cells.Item(14, 28).Value = -0.2043633904885378
$ws.Cells.Item(14, 25).Value = 0.1275780827634909
$ws.Cells.Item(14, 31).Value = -0.5849047489490333
$ws.Cells.Item(14, 27).Value = -0.01670081902098719
$ws.Cells.Item(14, 23).Value = 0.1240365846986169
$ws.Cells.Item(15, 30).Value = 0.0925667197466451
$ws.Cells.Item(15, 35).Value = -0.4334047671505248
$ws.Cells.Item(15, 32).Value = -1.40802832891157
$ws.Cells.Item(15, 33).Value = -1.407243743159736
$ws.Cells.Item(15, 31).Value = -0.2059746096811033
$ws.Cells.Item(15, 36).Value = -0.3096364143617802
$ws.Cells.Item(15, 27).Value = 0.09988477497939741
$ws.Cells.Item(15, 34).Value = -0.2500618974080826
$ws.Cells.Item(15, 28).Value = 0.08096036838765031
$ws.Cells.Item(15, 29).Value = 0.0924908932996793
$ws.Cells.Item(16, 36).Value = 0.2283024244226883
$ws.Cells.Item(16, 37).Value = 0.2048390592685578
$ws.Cells.Item(16, 31).Value = 0.02690490217465147
$ws.Cells.Item(16, 35).Value = -0.06335028919957075
$ws.Cells.Item(16, 39).Value = -0.1663214453978101
$ws.Cells.Item(16, 33).Value = -0.3318847693681293
$ws.Cells.Item(16, 38).Value = -0.0355780787674953
$ws.Cells.Item(16, 32).Value = -0.3373802571851825
$ws.Cells.Item(16, 40).Value = -0.1730430455425092
$ws.Cells.Item(16, 34).Value = 0.07829984441984905
$ws.Cells.Item(17, 35).Value = -0.1084535449743185
$ws.Cells.Item(17, 43).Value = 0.6502606143725664
$ws.Cells.Item(17, 40).Value = -0.4115424244148125
$ws.Cells.Item(17, 36).Value = -0.07323448430569535
$ws.Cells.Item(17, 38).Value = -0.1890623092888566
$ws.Cells.Item(17, 41).Value = 0.2152263639657814
$ws.Cells.Item(17, 37).Value = -0.09464543652764057
$ws.Cells.Item(17, 42).Value = 0.6050248749486009
$ws.Cells.Item(17, 34).Value = -0.1014535848389841
$ws.Cells.Item(17, 44).Value = 0.6376744206510576
$ws.Cells.Item(17, 39).Value = -0.3613518455741316
$ws.Cells.Item(18, 44).Value = -0.2058599286704377
$ws.Cells.Item(18, 47).Value = 0.1549171986535924
$ws.Cells.Item(18, 43).Value = -0.1691853834640433
$ws.Cells.Item(18, 48).Value = 0.1856341247700399
$ws.Cells.Item(18, 38).Value = -0.1350767357100935
$ws.Cells.Item(18, 46).Value = -0.05036452040672046
$ws.Cells.Item(18, 40).Value = -0.2817956528829213
$ws.Cells.Item(18, 39).Value = -0.2302498790168306
$ws.Cells.Item(18, 42).Value = -0.1865259660156937
$ws.Cells.Item(18, 41).Value = -0.2954722246111707
$ws.Cells.Item(18, 45).Value = -0.06071040501895997
$ws.Cells.Item(19, 48).Value = 0.0444032571666142
$ws.Cells.Item(19, 49).Value = 0.160714157635633
$ws.Cells.Item(19, 42).Value = -0.2229868532896306
$ws.Cells.Item(19, 46).Value = -0.2233431232791294
$ws.Cells.Item(19, 47).Value = -0.03486668218654065
$ws.Cells.Item(19, 50).Value = 0.01130372647704103
$ws.Cells.Item(19, 51).Value = -0.06391119588061711
$ws.Cells.Item(19, 52).Value = -0.09450306168263811
$ws.Cells.Item(19, 43).Value = -0.2399087565268632
$ws.Cells.Item(19, 44).Value = -0.2555246202002537
$ws.Cells.Item(19, 45).Value = -0.2428218159789997
$ws.Cells.Item(20, 46).Value = -0.214870093455366
$ws.Cells.Item(20, 50).Value = -0.1505299872523014
$ws.Cells.Item(20, 49).Value = -0.1040004763365077
$ws.Cells.Item(20, 51).Value = -0.2122873162357264
$ws.Cells.Item(20, 47).Value = -0.1571058138897907
$ws.Cells.Item(20, 48).Value = -0.1327251081615577
$ws.Cells.Item(20, 52).Value = -0.2824524929558314
$ws.Cells.Item(21, 50).Value = -0.1156228055083641
$ws.Cells.Item(21, 51).Value = -0.1284476411859137
$ws.Cells.Item(21, 52).Value = -0.1394689437024588
